$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "54.453.73"
$ws.Cells.Item(2, 5).Value = "  +0.42%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.287.08"
$ws.Cells.Item(3, 5).Value = "  +0.23%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.00%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "502.64"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +1.93%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "130.17"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +2.47%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -0.28%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.39%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.0959"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +2.15%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  +0.86%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  +5.00%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  +2.34%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "2.694.62"
$ws.Cells.Item(13, 5).Value = "  +0.34%  "

# Row 14
$ws.Cells.Item(14, 5).Value = "  +6.92%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "54.425.10"
$ws.Cells.Item(15, 5).Value = "  +0.68%  "

# Row 16
$ws.Cells.Item(16, 5).Value = "  +0.83%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "2.282.32"
$ws.Cells.Item(17, 5).Value = "  +0.76%  "

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "10.30"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +4.28%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  +3.07%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "304.89"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +1.16%  "

# Row 21
$ws.Cells.Item(21, 5).Value = "  +0.06%  "

# Row 22
$ws.Cells.Item(22, 5).Value = "  +0.20%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "61.94"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -2.76%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  -0.24%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  +2.24%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  +3.45%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "172.95"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +4.73%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  +2.47%  "

# Row 29
$ws.Cells.Item(29, 2).Value = "PEPE"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(29, 4).Value = "0.0₃0694"
$ws.Cells.Item(29, 5).Value = "  +2.08%  "

# Row 30
$ws.Cells.Item(30, 2).Value = "Aptos"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "5.99"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +2.04%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  +2.60%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  +1.91%  "

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.967"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +11.65%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  -0.27%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  +1.99%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "3.79"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +4.94%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  +0.34%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  +2.07%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "3.39"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +1.47%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "4.90"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +2.55%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "125.93"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +0.20%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  +4.08%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.0897"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +1.19%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.551"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +1.44%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "243.39"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +2.97%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  +0.31%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  +2.12%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  +0.92%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "16.51"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +2.00%  "

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "4.64"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -0.22%  "
